$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.725.67'
$ws.Range('E2').Value = '  +8.92%  '
$ws.Range('D3').Value = '2.689.38'
$ws.Range('E3').Value = '  +11.01%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'188.99"
$ws.Range('E5').Value = '  +13.61%  '
$ws.Range('D6').Value = "'589.85"
$ws.Range('E6').Value = '  +4.68%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').Value = "'0.542"
$ws.Range('E8').Value = '  +5.50%  '
$ws.Range('E9').Value = '  +14.99%  '
$ws.Range('D10').Value = '2.686.93'
$ws.Range('E10').Value = '  +10.98%  '
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = "'0.359"
$ws.Range('E12').Value = '  +7.32%  '
$ws.Range('D13').Value = "'4.76"
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').Value = '75.702.44'
$ws.Range('E14').Value = '  +9.16%  '
$ws.Range('D15').Value = '3.185.84'
$ws.Range('E15').Value = '  +10.98%  '
$ws.Range('D16').Value = "'0.0000190"
$ws.Range('E16').Value = '  +6.82%  '
$ws.Range('D17').Value = "'26.67"
$ws.Range('E17').Value = '  +11.20%  '
$ws.Range('D18').Value = '2.691.96'
$ws.Range('E18').Value = '  +10.80%  '
$ws.Range('D19').Value = "'9.42"
$ws.Range('E19').Value = '  +32.01%  '
$ws.Range('D20').Value = "'12.06"
$ws.Range('E20').Value = '  +11.55%  '
$ws.Range('D21').Value = "'375.05"
$ws.Range('E21').Value = '  +9.40%  '
$ws.Range('E22').Value = '  +17.92%  '
$ws.Range('D23').Value = "'4.08"
$ws.Range('E23').Value = '  +5.22%  '
$ws.Range('E24').Value = '  +4.03%  '
$ws.Range('D25').Value = "'0.999"
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').Value = "'70.42"
$ws.Range('E26').Value = '  +6.39%  '
$ws.Range('D27').Value = "'4.19"
$ws.Range('E27').Value = '  +9.70%  '
$ws.Range('D28').Value = "'9.46"
$ws.Range('E28').Value = '  +11.50%  '
$ws.Range('D29').Value = '2.831.08'
$ws.Range('E29').Value = '  +11.04%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('D31').Value = '0.0₃0956'
$ws.Range('E31').Value = '  +12.57%  '
$ws.Range('D32').Value = "'522.72"
$ws.Range('E32').Value = '  +15.08%  '
$ws.Range('E33').Value = '  +13.45%  '
$ws.Range('D34').Value = "'7.80"
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('D35').Value = "'1.76"
$ws.Range('D36').Value = "'0.999"
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = "'0.119"
$ws.Range('E37').Value = '  +8.50%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').Value = "'162.18"
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').Value = "'19.34"
$ws.Range('E39').Value = '  +6.23%  '
$ws.Range('D40').Value = "'19.38"
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('E42').Value = '  +14.50%  '
$ws.Range('D43').Value = "'170.87"
$ws.Range('E43').Value = '  +26.54%  '
$ws.Range('E44').Value = '  +12.31%  '
$ws.Range('E45').Value = '  +9.97%  '
$ws.Range('E46').Value = '  +10.67%  '
$ws.Range('D47').Value = "'2.39"
$ws.Range('E47').Value = '  +14.52%  '
$ws.Range('D48').Value = "'39.40"
$ws.Range('E48').Value = '  +4.22%  '
$ws.Range('D49').Value = "'0.0847"
$ws.Range('E49').Value = '  +16.86%  '
$ws.Range('E50').Value = '  +8.28%  '
$ws.Range('E51').Value = '  +10.64%  '
